# Update the cryptos price/volume table with freshly scraped figures.
# Note: several "Price" cells look numeric (e.g. "72.92") but must stay
# stored as text (as in the source data), so those assignments are
# prefixed with a leading apostrophe - the normal Excel trick that forces
# text entry instead of silently converting to a Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.111.00'
$ws.Range('E2').Value = '  +1.11%  '
$ws.Range('D3').Value = '2.152.16'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''253.59'
$ws.Range('E5').Value = '  +6.34%  '
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('D7').Value = '''72.92'
$ws.Range('E7').Value = '  +0.91%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '''0.578'
$ws.Range('E9').Value = '  -0.21%  '
$ws.Range('D10').Value = '''39.57'
$ws.Range('E10').Value = '  -1.15%  '
$ws.Range('D11').Value = '''0.0906'
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('E12').Value = '  +0.56%  '
$ws.Range('D13').Value = '''6.72'
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('D14').Value = '2.480.20'
$ws.Range('E14').Value = '  -0.63%  '
$ws.Range('E15').Value = '  -1.14%  '
$ws.Range('D16').Value = '2.122.71'
$ws.Range('E16').Value = '  -1.55%  '
$ws.Range('E17').Value = '  -2.07%  '
$ws.Range('D18').Value = '41.980.48'
$ws.Range('E18').Value = '  +1.16%  '
$ws.Range('D19').Value = '''0.0000102'
$ws.Range('E19').Value = '  -1.75%  '
$ws.Range('D20').Value = '''70.34'
$ws.Range('E20').Value = '  +0.32%  '
$ws.Range('E21').Value = '  +0.32%  '
$ws.Range('D22').Value = '''9.54'
$ws.Range('E22').Value = '  -2.51%  '
$ws.Range('D23').Value = '''225.47'
$ws.Range('E23').Value = '  -0.59%  '
$ws.Range('E24').Value = '  +4.65%  '
$ws.Range('E25').Value = '  -0.21%  '
$ws.Range('D26').Value = '''10.44'
$ws.Range('E26').Value = '  -2.86%  '
$ws.Range('E27').Value = '  +0.87%  '
$ws.Range('E28').Value = '  +2.76%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('D30').Value = '''36.80'
$ws.Range('E30').Value = '  +9.50%  '
$ws.Range('D31').Value = '''168.20'
$ws.Range('E31').Value = '  -1.64%  '
$ws.Range('D32').Value = '''19.85'
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('D33').Value = '''0.0798'
$ws.Range('E33').Value = '  +3.08%  '
$ws.Range('D34').Value = '''5.09'
$ws.Range('E34').Value = '  -3.00%  '
$ws.Range('E35').Value = '  -0.76%  '
$ws.Range('E36').Value = '  +1.36%  '
$ws.Range('E37').Value = '  -1.86%  '
$ws.Range('D38').Value = '''0.0330'
$ws.Range('E38').Value = '  +7.53%  '
$ws.Range('D39').Value = '''11.99'
$ws.Range('E39').Value = '  -1.58%  '
$ws.Range('E40').Value = '  -2.85%  '
$ws.Range('E41').Value = '  +3.50%  '
$ws.Range('D42').Value = '''5.13'
$ws.Range('E42').Value = '  -4.74%  '
$ws.Range('E43').Value = '  -1.05%  '
$ws.Range('D44').Value = '''99.78'
$ws.Range('E44').Value = '  +2.22%  '
$ws.Range('B45').Value = 'WOONetwork'
$ws.Range('C45').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D45').Value = '''0.463'
$ws.Range('E45').Value = '  +14.19%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '''8.22'
$ws.Range('E46').Value = '  -2.67%  '
$ws.Range('D47').Value = '''0.0962'
$ws.Range('E47').Value = '  -0.45%  '
$ws.Range('D48').Value = '''2.38'
$ws.Range('E48').Value = '  +8.64%  '
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('E50').Value = '  +0.47%  '
$ws.Range('E51').Value = '  +0.74%  '
